# Doing Updates for Financials
# Refresh the scraped NTZ yearly financial figures (income statement,
# balance sheet, and cash flow statement) with the latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement
$ws.Range("D8").Value = 504400
$ws.Range("E8").Value = 513000
$ws.Range("F8").Value = 548100
$ws.Range("G8").Value = 517700
$ws.Range("H8").Value = 503900
$ws.Range("I8").Value = 526000
$ws.Range("J8").Value = 545700
$ws.Range("D9").Value = 339000
$ws.Range("E9").Value = 336900
$ws.Range("F9").Value = 370900
$ws.Range("G9").Value = 373800
$ws.Range("H9").Value = 356000
$ws.Range("I9").Value = 352100
$ws.Range("J9").Value = 365800
$ws.Range("D10").Value = 165400
$ws.Range("E10").Value = 176100
$ws.Range("F10").Value = 177200
$ws.Range("G10").Value = 143900
$ws.Range("H10").Value = 147900
$ws.Range("I10").Value = 173900
$ws.Range("J10").Value = 179900
$ws.Range("D14").Value = 10400
$ws.Range("I14").Value = 1300
$ws.Range("J14").Value = 7600
$ws.Range("D15").Value = 4700
$ws.Range("E15").Value = 3900
$ws.Range("F15").Value = 4300
$ws.Range("G15").Value = 5100
$ws.Range("H15").Value = 7000
$ws.Range("D17").Value = 537100
$ws.Range("E17").Value = 513500
$ws.Range("F17").Value = 556600
$ws.Range("G17").Value = 559200
$ws.Range("H17").Value = 540200
$ws.Range("I17").Value = 546800
$ws.Range("J17").Value = 583900
$ws.Range("D18").Value = -32700
$ws.Range("F18").Value = -8600
$ws.Range("G18").Value = -41500
$ws.Range("H18").Value = -36300
$ws.Range("I18").Value = -20700
$ws.Range("J18").Value = -38200
$ws.Range("D20").Value = 6600
$ws.Range("E20").Value = 3300
$ws.Range("F20").Value = -4800
$ws.Range("G20").Value = -8700
$ws.Range("H20").Value = -33700
$ws.Range("J20").Value = 29100
$ws.Range("D21").Value = -11600
$ws.Range("E21").Value = 17500
$ws.Range("F21").Value = 2100
$ws.Range("G21").Value = -34200
$ws.Range("H21").Value = -51300
$ws.Range("I21").Value = -3400
$ws.Range("J21").Value = 12200
$ws.Range("D22").Value = 6800
$ws.Range("E22").Value = 5400
$ws.Range("F22").Value = 4500
$ws.Range("G22").Value = 3200
$ws.Range("H22").Value = 2100
$ws.Range("I22").Value = 1900
$ws.Range("J22").Value = 2100
$ws.Range("D23").Value = -32800
$ws.Range("E23").Value = -2500
$ws.Range("F23").Value = -17800
$ws.Range("G23").Value = -53300
$ws.Range("H23").Value = -72100
$ws.Range("I23").Value = -24500
$ws.Range("J23").Value = -11200
$ws.Range("D24").Value = 2900
$ws.Range("E24").Value = 4800
$ws.Range("F24").Value = 600
$ws.Range("G24").Value = 2000
$ws.Range("H24").Value = 4600
$ws.Range("I24").Value = 4700
$ws.Range("J24").Value = 19900
$ws.Range("D26").Value = -35700
$ws.Range("E26").Value = -7300
$ws.Range("F26").Value = -18500
$ws.Range("G26").Value = -55400
$ws.Range("H26").Value = -76700
$ws.Range("I26").Value = -29200
$ws.Range("J26").Value = -31100
$ws.Range("D27").Value = -35200
$ws.Range("E27").Value = -6800
$ws.Range("F27").Value = -18500
$ws.Range("G27").Value = -55400
$ws.Range("H27").Value = -76900
$ws.Range("I27").Value = -29300
$ws.Range("J27").Value = -24600
$ws.Range("D32").Value = -6600
$ws.Range("E32").Value = -3300
$ws.Range("F32").Value = 4800
$ws.Range("G32").Value = 8700
$ws.Range("H32").Value = 33700
$ws.Range("J32").Value = -29100
$ws.Range("D33").Value = -35200
$ws.Range("E33").Value = -6800
$ws.Range("F33").Value = -18500
$ws.Range("G33").Value = -55400
$ws.Range("H33").Value = -76900
$ws.Range("I33").Value = -29300
$ws.Range("J33").Value = -24600
$ws.Range("D35").Value = -35200
$ws.Range("E35").Value = -6800
$ws.Range("F35").Value = -18500
$ws.Range("G35").Value = -55400
$ws.Range("H35").Value = -76900
$ws.Range("I35").Value = -29300
$ws.Range("J35").Value = -24600
# Balance Sheet
$ws.Range("D41").Value = 61700
$ws.Range("E41").Value = 72900
$ws.Range("F41").Value = 58900
$ws.Range("G41").Value = 36900
$ws.Range("H41").Value = 68500
$ws.Range("I41").Value = 87200
$ws.Range("J41").Value = 105500
$ws.Range("D43").Value = 73500
$ws.Range("E43").Value = 88300
$ws.Range("F43").Value = 97700
$ws.Range("G43").Value = 128000
$ws.Range("H43").Value = 142900
$ws.Range("I43").Value = 276700
$ws.Range("J43").Value = 153400
$ws.Range("D44").Value = 90100
$ws.Range("E44").Value = 87900
$ws.Range("F44").Value = 88700
$ws.Range("G44").Value = 101200
$ws.Range("H44").Value = 88600
$ws.Range("I44").Value = 184600
$ws.Range("J44").Value = 104900
$ws.Range("D45").Value = 2200
$ws.Range("E45").Value = 3100
$ws.Range("F45").Value = 2400
$ws.Range("G45").Value = 2400
$ws.Range("H45").Value = 3200
$ws.Range("I45").Value = 3800
$ws.Range("J45").Value = 3400
$ws.Range("D46").Value = 227600
$ws.Range("E46").Value = 252300
$ws.Range("F46").Value = 247700
$ws.Range("G46").Value = 268500
$ws.Range("H46").Value = 303100
$ws.Range("I46").Value = 345000
$ws.Range("J46").Value = 367200
$ws.Range("E47").Value = 400
$ws.Range("F47").Value = 2500
$ws.Range("G47").Value = 3800
$ws.Range("H47").Value = 1600
$ws.Range("I47").Value = 1600
$ws.Range("J47").Value = 1600
$ws.Range("D48").Value = 121100
$ws.Range("E48").Value = 130100
$ws.Range("F48").Value = 135900
$ws.Range("G48").Value = 146700
$ws.Range("H48").Value = 161100
$ws.Range("I48").Value = 362400
$ws.Range("J48").Value = 197200
$ws.Range("D49").Value = 6200
$ws.Range("E49").Value = 4700
$ws.Range("F49").Value = 3800
$ws.Range("G49").Value = 4900
$ws.Range("H49").Value = 6200
$ws.Range("I49").Value = 5500
$ws.Range("J49").Value = 6000
$ws.Range("D52").Value = 1600
$ws.Range("E52").Value = 2100
$ws.Range("F52").Value = 2200
$ws.Range("G52").Value = 2500
$ws.Range("H52").Value = 1300
$ws.Range("I52").Value = 8800
$ws.Range("J52").Value = 1200
$ws.Range("D54").Value = 356500
$ws.Range("E54").Value = 389600
$ws.Range("F54").Value = 392000
$ws.Range("G54").Value = 426500
$ws.Range("H54").Value = 473400
$ws.Range("I54").Value = 534200
$ws.Range("J54").Value = 573300
$ws.Range("D57").Value = 85300
$ws.Range("E57").Value = 79800
$ws.Range("F57").Value = 66100
$ws.Range("G57").Value = 84400
$ws.Range("H57").Value = 75600
$ws.Range("I57").Value = 71100
$ws.Range("J57").Value = 71300
$ws.Range("D58").Value = 27500
$ws.Range("E58").Value = 33400
$ws.Range("F58").Value = 25100
$ws.Range("G58").Value = 26800
$ws.Range("H58").Value = 31800
$ws.Range("I58").Value = 34100
$ws.Range("J58").Value = 31300
$ws.Range("D59").Value = 52900
$ws.Range("E59").Value = 55300
$ws.Range("F59").Value = 49200
$ws.Range("G59").Value = 56800
$ws.Range("H59").Value = 47700
$ws.Range("I59").Value = 49200
$ws.Range("J59").Value = 35300
$ws.Range("D60").Value = 165700
$ws.Range("E60").Value = 168500
$ws.Range("F60").Value = 140400
$ws.Range("G60").Value = 168000
$ws.Range("H60").Value = 155100
$ws.Range("I60").Value = 149400
$ws.Range("J60").Value = 137800
$ws.Range("D61").Value = 23400
$ws.Range("E61").Value = 7100
$ws.Range("F61").Value = 17500
$ws.Range("G61").Value = 6900
$ws.Range("H61").Value = 4700
$ws.Range("I61").Value = 8200
$ws.Range("J61").Value = 12100
$ws.Range("D62").Value = 45700
$ws.Range("E62").Value = 42900
$ws.Range("F62").Value = 53900
$ws.Range("G62").Value = 56300
$ws.Range("H62").Value = 76200
$ws.Range("I62").Value = 58200
$ws.Range("J62").Value = 94400
$ws.Range("D66").Value = 237100
$ws.Range("E66").Value = 222400
$ws.Range("F66").Value = 215500
$ws.Range("G66").Value = 234600
$ws.Range("H66").Value = 239000
$ws.Range("I66").Value = 218700
$ws.Range("J66").Value = 224900
$ws.Range("D72").Value = 57900
$ws.Range("E72").Value = 105700
$ws.Range("F72").Value = 115000
$ws.Range("G72").Value = 120900
$ws.Range("H72").Value = 163400
$ws.Range("I72").Value = 244500
$ws.Range("J72").Value = 277500
$ws.Range("D76").Value = 119400
$ws.Range("E76").Value = 167200
$ws.Range("F76").Value = 176500
$ws.Range("G76").Value = 191900
$ws.Range("H76").Value = 234400
$ws.Range("I76").Value = 315500
$ws.Range("J76").Value = 348400
# Cash Flow Statement
$ws.Range("D81").Value = -35200
$ws.Range("E81").Value = -6800
$ws.Range("F81").Value = -18500
$ws.Range("G81").Value = -55400
$ws.Range("H81").Value = -76900
$ws.Range("I81").Value = -29300
$ws.Range("J81").Value = -24600
$ws.Range("D83").Value = 14400
$ws.Range("E83").Value = 14600
$ws.Range("F83").Value = 15400
$ws.Range("G83").Value = 16000
$ws.Range("H83").Value = 18600
$ws.Range("I83").Value = 19100
$ws.Range("J83").Value = 21300
$ws.Range("D89").Value = -5500
$ws.Range("E89").Value = 29200
$ws.Range("F89").Value = 9600
$ws.Range("G89").Value = -41700
$ws.Range("H89").Value = -2500
$ws.Range("I89").Value = -9200
$ws.Range("J89").Value = -22400
$ws.Range("D91").Value = -7500
$ws.Range("E91").Value = -6800
$ws.Range("F91").Value = -3900
$ws.Range("G91").Value = -7400
$ws.Range("H91").Value = -8000
$ws.Range("I91").Value = -8400
$ws.Range("J91").Value = -22100
$ws.Range("D94").Value = -13100
$ws.Range("E94").Value = -12000
$ws.Range("F94").Value = 1900
$ws.Range("G94").Value = 6500
$ws.Range("H94").Value = -9200
$ws.Range("I94").Value = -7300
$ws.Range("J94").Value = 31500
$ws.Range("D100").Value = 10400
$ws.Range("E100").Value = -2100
$ws.Range("F100").Value = 9000
$ws.Range("G100").Value = -2800
$ws.Range("H100").Value = -5800
$ws.Range("J100").Value = 26000
$ws.Range("D101").Value = -3000
$ws.Range("F101").Value = 1500
$ws.Range("G101").Value = 6400
$ws.Range("H101").Value = -1100
$ws.Range("J101").Value = 1800
$ws.Range("D102").Value = -11200
$ws.Range("E102").Value = 14000
$ws.Range("F102").Value = 22000
$ws.Range("G102").Value = -31600
$ws.Range("H102").Value = -18700
$ws.Range("I102").Value = -18300
$ws.Range("J102").Value = 37000
